$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scrape timestamp (column O) for every data row (2-398).
$newTimestamp = "2023-01-10 12:56:40"

for ($r = 2; $r -le 398; $r++) {
    $cell = $ws.Cells.Item($r, 15)  # column O
    $cell.Value = $newTimestamp
}

# A handful of productAriaLabel (column M) values had their
# "- Online kein Bestand" availability note added/removed.
$ws.Range("M23").Value = "Betty Bossi Frischback Buttergipfel IP-Suisse 2.60 Schweizer Franken"
$ws.Range("M181").Value = "Prix Garantie Roggenvollkornbrot - Online kein Bestand 1.80 Schweizer Franken"
$ws.Range("M244").Value = "Betty Bossi Vogelnestli 2x  100g - Online kein Bestand 4.60 Schweizer Franken"
$ws.Range("M377").Value = "Leisi Kuchenteig rund ausgewallt " + [char]0xD8 + "32cm glutenfrei 20% ab 2 Aktion 4.95 Schweizer Franken"
$ws.Range("M384").Value = "Leisi Bl" + [char]0xE4 + "tterteig glutenfrei rund ausgewallt " + [char]0xD8 + "32cm 20% ab 2 Aktion 4.95 Schweizer Franken"
